$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data gained one new weekly record. It is inserted as row 31
# (pushing every following record down by one row — the former row 111
# becomes row 112), then populated with its own values.
$ws.Rows.Item(31).Insert()

$ws.Cells.Item(31, 1).Value = 7
$ws.Cells.Item(31, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(31, 3).Value = "Ñuble"
$ws.Cells.Item(31, 4).Value = 44949
$ws.Cells.Item(31, 5).Value = 16
$ws.Cells.Item(31, 6).Value = 100112030
$ws.Cells.Item(31, 7).Value = "Poroto granado"
$ws.Cells.Item(31, 8).Value = "Sin especificar"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 50
$ws.Cells.Item(31, 11).Value = 40000
$ws.Cells.Item(31, 12).Value = 40000
$ws.Cells.Item(31, 13).Value = 40000
$ws.Cells.Item(31, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(31, 15).Value = "Región del Maule"
$ws.Cells.Item(31, 16).Value = 1600
$ws.Cells.Item(31, 17).Value = 25
$ws.Cells.Item(31, 18).Value = "Hortaliza"
